$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "73.246.43"
$ws.Range("E2").Value = "  -0.08%  "
$ws.Range("D3").Value = "3.985.38"
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'610.13"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.36%  "
$ws.Range("D6").Value = "'173.21"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +14.14%  "
$ws.Range("D7").Value = "'0.696"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.63%  "
$ws.Range("D8").Value = "'1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  +5.26%  "
$ws.Range("E10").Value = "  +9.08%  "
$ws.Range("D11").Value = "'56.82"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.53%  "
$ws.Range("D12").Value = "'0.0000338"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.63%  "
$ws.Range("D13").Value = "'11.77"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.37%  "
$ws.Range("D14").Value = "4.621.44"
$ws.Range("E14").Value = "  -1.66%  "
$ws.Range("D15").Value = "3.988.55"
$ws.Range("E15").Value = "  -1.57%  "
$ws.Range("D16").Value = "'14.34"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.11%  "
$ws.Range("B17").Value = "Chainlink"
$ws.Range("C17").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D17").Value = "'21.22"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.79%  "
$ws.Range("B18").Value = "Polygon"
$ws.Range("C18").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D18").Value = "'1.26"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.28%  "
$ws.Range("D19").Value = "73.197.23"
$ws.Range("E19").Value = "  -0.03%  "
$ws.Range("E20").Value = "  -0.81%  "
$ws.Range("D21").Value = "'466.25"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.84%  "
$ws.Range("D22").Value = "'4.84"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +6.77%  "
$ws.Range("D23").Value = "'97.75"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.38%  "
$ws.Range("E24").Value = "  -3.95%  "
$ws.Range("D25").Value = "'14.35"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.75%  "
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("D27").Value = "'11.40"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.09%  "
$ws.Range("E28").Value = "  -2.64%  "
$ws.Range("D29").Value = "'5.87"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.64%  "
$ws.Range("D30").Value = "'36.55"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.38%  "
$ws.Range("D31").Value = "'7.97"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.39%  "
$ws.Range("D32").Value = "'14.24"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.97%  "
$ws.Range("D33").Value = "'49.93"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.92%  "
$ws.Range("E34").Value = "  -2.86%  "
$ws.Range("E35").Value = "  +14.22%  "
$ws.Range("D36").Value = "'70.66"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.51%  "
$ws.Range("D37").Value = "'643.53"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -6.98%  "
$ws.Range("D38").Value = "'0.436"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.48%  "
$ws.Range("B39").Value = "ThetaToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D39").Value = "'3.43"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.63%  "
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").Value = "'0.149"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.87%  "
$ws.Range("B41").Value = "dogwifhat"
$ws.Range("C41").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D41").Value = "'3.38"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +49.11%  "
$ws.Range("B42").Value = "Dai"
$ws.Range("C42").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D42").Value = "'1.00"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.03%  "
$ws.Range("E43").Value = "  +0.16%  "
$ws.Range("D44").Value = "'0.0488"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.56%  "
$ws.Range("D45").Value = "'10.57"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -7.65%  "
$ws.Range("D46").Value = "'0.151"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.28%  "
$ws.Range("D47").Value = "'3.03"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -8.98%  "
$ws.Range("B48").Value = "Fetch.AI"
$ws.Range("C48").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D48").Value = "'2.70"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.37%  "
$ws.Range("B49").Value = "FLOKI"
$ws.Range("C49").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D49").Value = "'0.000298"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +9.27%  "
$ws.Range("D50").Value = "'3.44"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.59%  "
$ws.Range("D51").Value = "2.829.60"
$ws.Range("E51").Value = "  +1.16%  "
